$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Holly added "S.GISH" to the harvester list in bioSamples, which fixes the
# harvester column (B) here in rnaSamples -- update every data row (2-29).
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 2).Value = "S.GISH"
}

# Column B widened slightly to fit the new text.
$ws.Columns.Item(2).ColumnWidth = 8

# Leave the selection on the harvester column, matching the post-edit state.
$ws.Range("B:B").Select()
